# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change cell B11 on the "Rules" sheet from the shared string "R40" to the
# literal text "1" (a string, not the number 1), keeping the cell's
# existing style (border/format) untouched.
#
# A plain  Range.Value = "1"  assignment gets auto-coerced to the number 1
# by the workbook's type inference, and prefixing with an apostrophe
# ("'1") forces text but also flips the cell format's quote-prefix flag,
# which creates a brand-new style entry (so the cell's style index would
# change). Neither matches the intended edit, which keeps the same style.
#
# Instead: compute the text "1" via a formula in a scratch cell (TEXT()
# always returns a genuine string), copy just that *value* (not the
# formatting) onto B11, then remove the scratch cell again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z100")
$target = $ws.Range("B11")

$scratch.Formula = '=TEXT(1,"0")'
$scratch.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues: value only, formatting untouched
$scratch.Clear()
